$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.438.34"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "2.240.55"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "1.02"
$ws.Range("E4").Value = "  +1.29%  "
$ws.Range("D5").Value = "307.37"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "94.33"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "34.70"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").Value = "0.0803"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "7.20"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D14").Value = "2.270.13"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "13.54"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "44.097.53"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "11.90"
$ws.Range("E20").Value = "  -2.65%  "
$ws.Range("D21").Value = "65.49"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "237.49"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").Value = "2.95"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").Value = "38.14"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").Value = "9.77"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").Value = "5.95"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "19.92"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "153.01"
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("D32").Value = "0.0794"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "2.64"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").Value = "3.08"
$ws.Range("E34").Value = "  -3.71%  "
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").Value = "14.94"
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").Value = "3.74"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "1.790.59"
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").Value = "78.50"
$ws.Range("E45").Value = "  -8.39%  "
$ws.Range("D46").Value = "70.14"
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "98.50"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "4.88"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  +4.84%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  +0.21%  "
